$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated financial figures (row 3-35) ---
$ws.Range("B3").Value = 513983
$ws.Range("C3").Value = 469822
$ws.Range("D3").Value = 386064
$ws.Range("E3").Value = 280522
$ws.Range("B4").Value = 513983
$ws.Range("C4").Value = 469822
$ws.Range("D4").Value = 386064
$ws.Range("E4").Value = 280522
$ws.Range("B6").Value = 288831
$ws.Range("C6").Value = 272344
$ws.Range("D6").Value = 233307
$ws.Range("E6").Value = 165536
$ws.Range("B7").Value = 225152
$ws.Range("C7").Value = 197478
$ws.Range("D7").Value = 152757
$ws.Range("E7").Value = 114986
$ws.Range("B8").Value = 501735
$ws.Range("C8").Value = 444943
$ws.Range("D8").Value = 363165
$ws.Range("E8").Value = 265981
$ws.Range("B9").Value = 138428
$ws.Range("C9").Value = 116485
$ws.Range("D9").Value = 87193
$ws.Range("E9").Value = 64313
$ws.Range("B10").Value = 73213
$ws.Range("C10").Value = 56052
$ws.Range("D10").Value = 42740
$ws.Range("E10").Value = 35931
$ws.Range("B11").Value = 16393
$ws.Range("C11").Value = 11012
$ws.Range("D11").Value = 8432
$ws.Range("E11").Value = 565
$ws.Range("B12").Value = -2367
$ws.Range("C12").Value = -1809
$ws.Range("D12").Value = -1647
$ws.Range("E12").Value = "-"
$ws.Range("B14").Value = -12763
$ws.Range("C14").Value = -9141
$ws.Range("D14").Value = -6860
$ws.Range("E14").Value = -364
$ws.Range("B15").Value = 12248
$ws.Range("C15").Value = 24879
$ws.Range("D15").Value = 22899
$ws.Range("E15").Value = 14541
$ws.Range("C16").Value = 10110
$ws.Range("D16").Value = -224
$ws.Range("E16").Value = -557
$ws.Range("B17").Value = -15926
$ws.Range("C17").Value = 14707
$ws.Range("D17").Value = 2394
$ws.Range("E17").Value = "-"
$ws.Range("B18").Value = 34113
$ws.Range("C18").Value = -27983
$ws.Range("D18").Value = -3689
$ws.Range("E18").Value = -8
$ws.Range("B19").Value = -5939
$ws.Range("C19").Value = 38155
$ws.Range("D19").Value = 24194
$ws.Range("E19").Value = 13976
$ws.Range("B20").Value = -3217
$ws.Range("C20").Value = 4791
$ws.Range("D20").Value = 2863
$ws.Range("E20").Value = 2374
$ws.Range("B21").Value = -2722
$ws.Range("C21").Value = 33364
$ws.Range("D21").Value = 21331
$ws.Range("E21").Value = 11602
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 16
$ws.Range("E23").Value = -14
$ws.Range("B25").Value = -2722
$ws.Range("C25").Value = 33364
$ws.Range("D25").Value = 21331
$ws.Range("E25").Value = 11588
$ws.Range("B27").Value = -2722
$ws.Range("C27").Value = 33364
$ws.Range("D27").Value = 21331
$ws.Range("E27").Value = 11588
$ws.Range("B29").Value = -2722
$ws.Range("C29").Value = 33364
$ws.Range("D29").Value = 21331
$ws.Range("E29").Value = 11588
$ws.Range("B30").Value = 29.03
$ws.Range("C30").Value = 4.96
$ws.Range("D30").Value = 17.18
$ws.Range("E30").Value = "-"
$ws.Range("B31").Value = -2751.03
$ws.Range("C31").Value = 33359.04
$ws.Range("D31").Value = 21313.82
$ws.Range("E31").Value = 11588
$ws.Range("B32").Value = 10189
$ws.Range("C32").Value = 10296
$ws.Range("D32").Value = 10198
$ws.Range("E32").Value = 504
$ws.Range("B33").Value = -0.27
$ws.Range("C33").Value = 3.24
$ws.Range("D33").Value = 2.09
$ws.Range("E33").Value = 22.99
$ws.Range("B35").Value = 0.613
$ws.Range("C35").Value = 1.42
$ws.Range("D35").Value = 1.34
$ws.Range("E35").Value = 22.99

# --- Row shading updates: Ventas netas totales row gets a slightly darker blue,
# Dilucion de las ganancias por accion basicas row gets the light-blue band it was missing ---
$ws.Range("A3:E3").Interior.Color = 15919583
$ws.Range("A35:E35").Interior.Color = 16446701

# --- Restore the active selection to A3 ---
$null = $ws.Range("A3").Select()
